# Refresh the crypto price/volume table (mirrors the "Updated cryptos list" GitHub Action commit).
# For each changed cell we just overwrite the inline-string value. A handful of D-column prices
# (e.g. "0.9990", "92.00") are valid numeric literals, so plain assignment would make Excel
# auto-convert them to numbers and silently drop trailing zeros / text formatting. To avoid that we
# enter those with a leading apostrophe (forces text) and then reset .Style to "Normal" so no extra
# number-format/quote-prefix style sticks to the cell - the sheet stays byte-for-byte a plain inline string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.771.83'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '1.868.88'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('D4').Value = '''0.9990'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''0.7295'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.70%  '
$ws.Range('D6').Value = '''240.78'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.42%  '
$ws.Range('D7').Value = '''0.9993'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '''0.3122'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = '''0.07092'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.31%  '
$ws.Range('D10').Value = '''24.27'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('E11').Value = '  -2.16%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.884.92'
$ws.Range('E12').Value = '  +1.80%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').Value = '''0.7401'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.62%  '
$ws.Range('D14').Value = '''5.313'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').Value = '''92.00'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('D16').Value = '29.772.65'
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D17').Value = '''6.004'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '''247.64'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.04%  '
$ws.Range('D19').Value = '''13.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.64%  '
$ws.Range('D20').Value = '''0.000007791'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = '''0.9994'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.132.95'
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('D23').Value = '''0.9993'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '''7.733'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.52%  '
$ws.Range('D25').Value = '''0.1539'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.08%  '
$ws.Range('D26').Value = '''9.170'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.37%  '
$ws.Range('D27').Value = '''162.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').Value = '''18.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.51%  '
$ws.Range('D29').Value = '''2.004'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.84%  '
$ws.Range('D30').Value = '''1.440'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.18%  '
$ws.Range('D31').Value = '''4.520'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.45%  '
$ws.Range('D32').Value = '''1.520'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.57%  '
$ws.Range('D33').Value = '''4.162'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.45%  '
$ws.Range('D34').Value = '''0.05270'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.93%  '
$ws.Range('D35').Value = '''1.231'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').Value = '''0.7435'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('D38').Value = '''2.688'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('E39').Value = '  -1.38%  '
$ws.Range('E40').Value = '  -0.64%  '
$ws.Range('D41').Value = '''0.4448'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.53%  '
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('D43').Value = '''0.8644'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.54%  '
$ws.Range('D44').Value = '''71.15'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.36%  '
$ws.Range('D45').Value = '1.042.33'
$ws.Range('E45').Value = '  -5.51%  '
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').Value = '''103.69'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '''7.440'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.11%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '''1.811'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''9.478'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.022.18'
$ws.Range('E51').Value = '  +0.44%  '
